{"js": "// The questionnaire table has five bold section headings of the form\n// \"Responses of end-users and IT experts on PIECES Software Evaluation\n// ...\" (one plain, the others prefixed \"B./D./E./F.\"). Per the commit,\n// these should read \"Responses of end-users on PIECES Software\n// Evaluation ...\" i.e. drop \" and IT experts\". The search phrase below\n// is specific to those five headings, so it does not touch the\n// unrelated cover-page sentence (\"... assessed and evaluated by the\n// end-users and IT experts.\").\nconst body = context.document.body;\n\nconst results = body.search(\"end-users and IT experts on PIECES Software Evaluation\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"end-users on PIECES Software Evaluation\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The questionnaire table has five bold section headings of the form\n# \"Responses of end-users and IT experts on PIECES Software Evaluation\n# ...\" (and one prefixed \"B./D./E./F.\"). Per the commit, these should\n# read \"Responses of end-users on PIECES Software Evaluation ...\" i.e.\n# drop \" and IT experts\". The unrelated cover-page sentence\n# (\"... assessed and evaluated by the end-users and IT experts.\") must\n# stay untouched, so we scope the Find/Replace to the table that holds\n# the headings and use a search phrase specific to those headings.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$find = $tbl.Range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"end-users and IT experts on PIECES Software Evaluation\"\n$find.Replacement.Text = \"end-users on PIECES Software Evaluation\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $false\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2)\n"}
